$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.176.21"
$ws.Range("E2").Value = "  +4.26%  "

# Row 3
$ws.Range("D3").Value = "2.733.55"
$ws.Range("E3").Value = "  +2.94%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.48"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.44"
$ws.Range("E6").Value = "  +9.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +4.58%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.996"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").Value = "2.758.88"
$ws.Range("E9").Value = "  +3.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("E10").Value = "  +3.08%  "

# Row 11
$ws.Range("E11").Value = "  +2.91%  "

# Row 12
$ws.Range("E12").Value = "  +3.56%  "

# Row 13
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("D14").Value = "3.226.14"
$ws.Range("E14").Value = "  +3.12%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.36"
$ws.Range("E15").Value = "  +2.58%  "

# Row 16
$ws.Range("D16").Value = "63.771.44"
$ws.Range("E16").Value = "  +3.65%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000156"
$ws.Range("E17").Value = "  +6.16%  "

# Row 18
$ws.Range("D18").Value = "2.752.13"
$ws.Range("E18").Value = "  +3.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.11"
$ws.Range("E19").Value = "  +4.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.96"
$ws.Range("E20").Value = "  +3.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.08"
$ws.Range("E21").Value = "  +2.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.04"
$ws.Range("E22").Value = "  +1.40%  "

# Row 23
$ws.Range("E23").Value = "  +2.97%  "

# Row 24
$ws.Range("E24").Value = "  -0.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.87"
$ws.Range("E25").Value = "  +4.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.173"
$ws.Range("E26").Value = "  +5.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.66"
$ws.Range("E27").Value = "  +0.72%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.28%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0923"
$ws.Range("E29").Value = "  +12.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +1.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.26"
$ws.Range("E31").Value = "  +5.55%  "

# Row 32
$ws.Range("E32").Value = "  +13.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "173.65"
$ws.Range("E33").Value = "  +4.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.69"
$ws.Range("E35").Value = "  +3.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  +5.04%  "

# Row 37
$ws.Range("E37").Value = "  +9.07%  "

# Row 38
$ws.Range("E38").Value = "  +5.62%  "

# Row 39
$ws.Range("E39").Value = "  +10.70%  "

# Row 40
$ws.Range("E40").Value = "  +4.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.28"
$ws.Range("E41").Value = "  +17.30%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "339.57"
$ws.Range("E42").Value = "  -1.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.47"
$ws.Range("E43").Value = "  +2.86%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.91"
$ws.Range("E44").Value = "  +6.54%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.44"
$ws.Range("E45").Value = "  +5.80%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0603"
$ws.Range("E46").Value = "  +3.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.647"
$ws.Range("E47").Value = "  +2.77%  "

# Row 48
$ws.Range("E48").Value = "  +3.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.99"
$ws.Range("E49").Value = "  +2.02%  "

# Row 50
$ws.Range("E50").Value = "  +3.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  +0.00%  "
